# Update the "Förändrad" (Changed) date column (C) for rows 2-16.
# The stored serial date value moves from 46074 to 46075 (i.e. +1 day).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 16; $row++) {
    $cell = $ws.Cells.Item($row, 3)  # Column C
    if ($cell.Value2 -eq 46074) {
        $cell.Value = 46075
    }
}
